# Swap the presentation's theme colours from the "Integral" palette to the
# default "Office Theme" palette (theme1.xml governs the slide master that
# every slide/layout ultimately inherits from).
#
# PowerPoint's COM model doesn't expose a "replace whole theme part" call;
# the supported way to edit a theme's look is through ThemeColorScheme
# (a.k.a. Colors()) entries, one RGB colour at a time - in the fixed order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
#
# .RGB takes a VBA-style RGB() long: red + (green*256) + (blue*65536), i.e.
# the hex digits of an "RRGGBB" colour have to be written back-to-front as
# 0xBBGGRR.

$p   = $ppt.ActivePresentation
$cs  = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Colors(1).RGB  = 0x000000   # dk1      RRGGBB=000000
$cs.Colors(2).RGB  = 0xFFFFFF   # lt1      RRGGBB=FFFFFF
$cs.Colors(3).RGB  = 0x6A5444   # dk2      RRGGBB=44546A
$cs.Colors(4).RGB  = 0xE6E6E7   # lt2      RRGGBB=E7E6E6
$cs.Colors(5).RGB  = 0xD59B5B   # accent1  RRGGBB=5B9BD5
$cs.Colors(6).RGB  = 0x317DED   # accent2  RRGGBB=ED7D31
$cs.Colors(7).RGB  = 0xA5A5A5   # accent3  RRGGBB=A5A5A5
$cs.Colors(8).RGB  = 0x00C0FF   # accent4  RRGGBB=FFC000
$cs.Colors(9).RGB  = 0xC47244   # accent5  RRGGBB=4472C4
$cs.Colors(10).RGB = 0x47AD70   # accent6  RRGGBB=70AD47
$cs.Colors(11).RGB = 0xC16305   # hlink    RRGGBB=0563C1
$cs.Colors(12).RGB = 0x724F95   # folHlink RRGGBB=954F72
